# Update the RGB backtesting results on "sheet0" to reflect the
# re-run / refactored (object-oriented) model output.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet0")

$ws.Range("B2").Value = "RED"
$ws.Range("C2").Value = "RED"
$ws.Range("D2").Value = "RED"
$ws.Range("E2").Value = "RED"

$ws.Range("B3").Value = "RED"
$ws.Range("C3").Value = "RED"
$ws.Range("D3").Value = "GREEN"
$ws.Range("E3").Value = "GREEN"

$ws.Range("B4").Value = "RED"
$ws.Range("C4").Value = "RED"
$ws.Range("D4").Value = "RED"
$ws.Range("E4").Value = "RED"

$ws.Range("B5").Value = "RED"
$ws.Range("C5").Value = "RED"
$ws.Range("D5").Value = "RED"
$ws.Range("E5").Value = "RED"

$ws.Range("B6").Value = "AMBER"
$ws.Range("C6").Value = "RED"
$ws.Range("D6").Value = "RED"
$ws.Range("E6").Value = "AMBER"

$ws.Range("B7").Value = "RED"
$ws.Range("C7").Value = "AMBER"
$ws.Range("D7").Value = "AMBER"
$ws.Range("E7").Value = "GREEN"

$ws.Range("B8").Value = "GREEN"
$ws.Range("C8").Value = "GREEN"
$ws.Range("D8").Value = "GREEN"
$ws.Range("E8").Value = "GREEN"

$ws.Range("B9").Value = "RED"
$ws.Range("C9").Value = "RED"
$ws.Range("D9").Value = "RED"
$ws.Range("E9").Value = "RED"

$ws.Range("B10").Value = "RED"
$ws.Range("C10").Value = "RED"
$ws.Range("D10").Value = "RED"
$ws.Range("E10").Value = "RED"

$ws.Range("B11").Value = "RED"
$ws.Range("C11").Value = "AMBER"
$ws.Range("D11").Value = "GREEN"
$ws.Range("E11").Value = "GREEN"

$ws.Range("B12").Value = "GREEN"
$ws.Range("C12").Value = "GREEN"
$ws.Range("D12").Value = "GREEN"
$ws.Range("E12").Value = "GREEN"

$ws.Range("B13").Value = "RED"
$ws.Range("C13").Value = "RED"
$ws.Range("D13").Value = "AMBER"
$ws.Range("E13").Value = "AMBER"

$ws.Range("B14").Value = "RED"
$ws.Range("C14").Value = "RED"
$ws.Range("D14").Value = "GREEN"
$ws.Range("E14").Value = "GREEN"

$ws.Range("B15").Value = "RED"
$ws.Range("C15").Value = "RED"
$ws.Range("D15").Value = "AMBER"
$ws.Range("E15").Value = "GREEN"

$ws.Range("B16").Value = "RED"
$ws.Range("C16").Value = "GREEN"
$ws.Range("D16").Value = "GREEN"
$ws.Range("E16").Value = "GREEN"

$ws.Range("B17").Value = "AMBER"
$ws.Range("C17").Value = "GREEN"
$ws.Range("D17").Value = "GREEN"
$ws.Range("E17").Value = "GREEN"

$ws.Range("B18").Value = "RED"
$ws.Range("C18").Value = "GREEN"
$ws.Range("D18").Value = "GREEN"
$ws.Range("E18").Value = "GREEN"

$ws.Range("B19").Value = "RED"
$ws.Range("C19").Value = "RED"
$ws.Range("D19").Value = "RED"
$ws.Range("E19").Value = "AMBER"

$ws.Range("B20").Value = "RED"
$ws.Range("C20").Value = "RED"
$ws.Range("D20").Value = "RED"
$ws.Range("E20").Value = "RED"

$ws.Range("B21").Value = "RED"
$ws.Range("C21").Value = "RED"
$ws.Range("D21").Value = "RED"
$ws.Range("E21").Value = "RED"
